# Update gh-pages output values (F column "想去人数") for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 779
$ws1.Range("F11").Value = 479
$ws1.Range("F12").Value = 1432
$ws1.Range("F23").Value = 6169
$ws1.Range("F28").Value = 14924
$ws1.Range("F33").Value = 10877
$ws1.Range("F34").Value = 692
$ws1.Range("F35").Value = 4250
$ws1.Range("F36").Value = 192

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 779
$ws4.Range("F11").Value = 479
$ws4.Range("F12").Value = 1432
$ws4.Range("F26").Value = 6169
$ws4.Range("F31").Value = 14924
$ws4.Range("F36").Value = 10877
$ws4.Range("F37").Value = 692
$ws4.Range("F38").Value = 4250
$ws4.Range("F39").Value = 192
